$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").Value = "'" + '63.360.62'
$ws.Range("D2").Style = $style
$ws.Range("E2").Value = '  +0.28%  '

$style = $ws.Range("D3").Style
$ws.Range("D3").Value = "'" + '2.670.12'
$ws.Range("D3").Style = $style
$ws.Range("E3").Value = '  +3.60%  '

$ws.Range("E4").Value = '  -0.05%  '

$style = $ws.Range("D5").Style
$ws.Range("D5").Value = "'" + '611.48'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +3.17%  '

$style = $ws.Range("D6").Style
$ws.Range("D6").Value = "'" + '143.29'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  -0.94%  '

$ws.Range("E7").Value = '  -0.05%  '

$style = $ws.Range("D8").Style
$ws.Range("D8").Value = "'" + '0.586'
$ws.Range("D8").Style = $style
$ws.Range("E8").Value = '  -0.75%  '

$style = $ws.Range("D9").Style
$ws.Range("D9").Value = "'" + '2.669.03'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  +3.60%  '

$ws.Range("E10").Value = '  +0.24%  '

$style = $ws.Range("D12").Style
$ws.Range("D12").Value = "'" + '0.152'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  +0.33%  '

$ws.Range("E13").Value = '  +3.12%  '

$style = $ws.Range("D14").Style
$ws.Range("D14").Value = "'" + '27.34'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value = '  +0.44%  '

$style = $ws.Range("D15").Style
$ws.Range("D15").Value = "'" + '3.151.64'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  +3.62%  '

$style = $ws.Range("D16").Style
$ws.Range("D16").Value = "'" + '63.194.11'
$ws.Range("D16").Style = $style
$ws.Range("E16").Value = '  +0.17%  '

$ws.Range("E17").Value = '  -0.96%  '

$style = $ws.Range("D18").Style
$ws.Range("D18").Value = "'" + '2.689.09'
$ws.Range("D18").Style = $style
$ws.Range("E18").Value = '  +4.05%  '

$style = $ws.Range("D19").Style
$ws.Range("D19").Value = "'" + '11.42'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  +2.91%  '

$style = $ws.Range("D20").Style
$ws.Range("D20").Value = "'" + '341.38'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  -0.02%  '

$ws.Range("E21").Value = '  +1.49%  '

$ws.Range("E22").Value = '  +3.41%  '

$ws.Range("E23").Value = '  -0.08%  '

$style = $ws.Range("D24").Style
$ws.Range("D24").Value = "'" + '67.20'
$ws.Range("D24").Style = $style
$ws.Range("E24").Value = '  -0.94%  '

$style = $ws.Range("D25").Style
$ws.Range("D25").Value = "'" + '1.64'
$ws.Range("D25").Style = $style
$ws.Range("E25").Value = '  +1.11%  '

$ws.Range("E26").Value = '  -3.96%  '

$ws.Range("E27").Value = '  +4.74%  '

$ws.Range("E28").Value = '  -1.10%  '

$style = $ws.Range("D29").Style
$ws.Range("D29").Value = "'" + '543.58'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value = '  +15.54%  '

$ws.Range("E30").Value = '  +0.06%  '

$style = $ws.Range("D31").Style
$ws.Range("D31").Value = "'" + '7.90'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  -0.26%  '

$style = $ws.Range("D32").Style
$ws.Range("D32").Value = "'" + '2.06'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  +5.28%  '

$ws.Range("E33").Value = '  +6.95%  '

$ws.Range("E34").Value = '  +0.36%  '

$style = $ws.Range("D35").Style
$ws.Range("D35").Value = "'" + '172.43'
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  -2.38%  '

$style = $ws.Range("D36").Style
$ws.Range("D36").Value = "'" + '5.15'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  +12.79%  '

$ws.Range("E37").Value = '  -0.05%  '

$style = $ws.Range("D39").Style
$ws.Range("D39").Value = "'" + '19.18'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  +1.62%  '

$ws.Range("E40").Value = '  +9.61%  '

$style = $ws.Range("D41").Style
$ws.Range("D41").Value = "'" + '178.19'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  +12.25%  '

$ws.Range("E42").Value = '  +0.02%  '

$ws.Range("E43").Value = '  +0.90%  '

$ws.Range("E44").Value = '  +4.00%  '

$style = $ws.Range("D45").Style
$ws.Range("D45").Value = "'" + '0.0574'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value = '  +6.39%  '

$style = $ws.Range("D46").Style
$ws.Range("D46").Value = "'" + '0.635'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +0.00%  '

$style = $ws.Range("D47").Style
$ws.Range("D47").Value = "'" + '0.0963'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value = '  -0.13%  '

$ws.Range("E48").Value = '  +1.00%  '

$style = $ws.Range("D49").Style
$ws.Range("D49").Value = "'" + '18.70'
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  +2.98%  '

$style = $ws.Range("D50").Style
$ws.Range("D50").Value = "'" + '1.74'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  +3.38%  '

$ws.Range("E51").Value = '  -0.79%  '
